$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 11
$ws.Range("C4").Value = 5

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 5

# Row 18
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 11
$ws.Range("E18").Value = 5
